$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.9965699487229559
$ws.Range("D2").Value = 0.002208650319530431
$ws.Range("E2").Value = 0.9990623130852649

$ws.Range("C3").Value = -0.3265723221688601
$ws.Range("D3").Value = 0.4824760960532301
$ws.Range("E3").Value = -0.5873498981461435

$ws.Range("C4").Value = -0.05734928802674788
$ws.Range("D4").Value = 0.1249142439989834
$ws.Range("E4").Value = -0.0161497335617524

$ws.Range("C5").Value = -0.5359341722226971
$ws.Range("D5").Value = 0.5554471686140796
$ws.Range("E5").Value = -0.5637216189927574

$ws.Range("C6").Value = -10.95209495043962
$ws.Range("D6").Value = 0.8934283804811416
$ws.Range("E6").Value = -0.4776221073808635

$ws.Range("C7").Value = 0.9051383052607735
$ws.Range("D7").Value = 0.01443863690947908
$ws.Range("E7").Value = 0.9598009007022293

$ws.Range("C8").Value = -1.098025891347283
$ws.Range("D8").Value = 0.6698770340678948
$ws.Range("E8").Value = -0.3662337739607429

$ws.Range("C9").Value = 0.4242379068891781
$ws.Range("D9").Value = 0.01157189805751288
$ws.Range("E9").Value = 0.9750035643357913
